$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 49.400308
$ws.Range("H2").Value = 148.200924
$ws.Range("I2").Value = 0.3028101582105581
$ws.Range("J2").Value = 0.3028101582105581
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.111991333333334
$ws.Range("N2").Value = 9.335974
$ws.Range("O2").Value = 0.1375443917436657
$ws.Range("P2").Value = 0.1375443917436657
$ws.Range("Q2").Value = 153.7333303599973
$ws.Range("R2").Value = 1383.599973239976
$ws.Range("S2").Value = 0.04164983902487439
$ws.Range("T2").Value = 0.04164983902487438

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 49.400308
$ws.Range("H3").Value = 148.200924
$ws.Range("I3").Value = 0.3028101582105581
$ws.Range("J3").Value = 0.3028101582105581
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.729556333333333
$ws.Range("N3").Value = 14.188669
$ws.Range("O3").Value = 0.2090378408570124
$ws.Range("P3").Value = 0.2090378408570124
$ws.Range("Q3").Value = 233.6415395700173
$ws.Range("R3").Value = 2102.773856130156
$ws.Range("S3").Value = 0.06329878166190538
$ws.Range("T3").Value = 0.06329878166190538

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 49.400308
$ws.Range("H4").Value = 148.200924
$ws.Range("I4").Value = 0.3028101582105581
$ws.Range("J4").Value = 0.3028101582105581
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 14.783812
$ws.Range("N4").Value = 44.351436
$ws.Range("O4").Value = 0.6534177673993219
$ws.Range("P4").Value = 0.6534177673993219
$ws.Range("Q4").Value = 730.3248662140959
$ws.Range("R4").Value = 6572.923795926863
$ws.Range("S4").Value = 0.1978615375237783
$ws.Range("T4").Value = 0.1978615375237783

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 69.564149
$ws.Range("H5").Value = 208.692447
$ws.Range("I5").Value = 0.4264088994034781
$ws.Range("J5").Value = 0.4264088994034782
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.111991333333334
$ws.Range("N5").Value = 9.335974
$ws.Range("O5").Value = 0.1375443917436657
$ws.Range("P5").Value = 0.1375443917436657
$ws.Range("Q5").Value = 216.4830287987087
$ws.Range("R5").Value = 1948.347259188378
$ws.Range("S5").Value = 0.05865015270253734
$ws.Range("T5").Value = 0.05865015270253733

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 69.564149
$ws.Range("H6").Value = 208.692447
$ws.Range("I6").Value = 0.4264088994034781
$ws.Range("J6").Value = 0.4264088994034782
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 4.729556333333333
$ws.Range("N6").Value = 14.188669
$ws.Range("O6").Value = 0.2090378408570124
$ws.Range("P6").Value = 0.2090378408570124
$ws.Range("Q6").Value = 329.0075614758936
$ws.Range("R6").Value = 2961.068053283043
$ws.Range("S6").Value = 0.08913559565351806
$ws.Range("T6").Value = 0.08913559565351807

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 69.564149
$ws.Range("H7").Value = 208.692447
$ws.Range("I7").Value = 0.4264088994034781
$ws.Range("J7").Value = 0.4264088994034782
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 14.783812
$ws.Range("N7").Value = 44.351436
$ws.Range("O7").Value = 0.6534177673993219
$ws.Range("P7").Value = 0.6534177673993219
$ws.Range("Q7").Value = 1028.423300755988
$ws.Range("R7").Value = 9255.809706803893
$ws.Range("S7").Value = 0.2786231510474227
$ws.Range("T7").Value = 0.2786231510474227

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 44.17507666666666
$ws.Range("H8").Value = 132.52523
$ws.Range("I8").Value = 0.2707809423859638
$ws.Range("J8").Value = 0.2707809423859638
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.111991333333334
$ws.Range("N8").Value = 9.335974
$ws.Range("O8").Value = 0.1375443917436657
$ws.Range("P8").Value = 0.1375443917436657
$ws.Range("Q8").Value = 137.4724557360022
$ws.Range("R8").Value = 1237.25210162402
$ws.Range("S8").Value = 0.03724440001625397
$ws.Range("T8").Value = 0.03724440001625397

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 44.17507666666666
$ws.Range("H9").Value = 132.52523
$ws.Range("I9").Value = 0.2707809423859638
$ws.Range("J9").Value = 0.2707809423859638
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.729556333333333
$ws.Range("N9").Value = 14.188669
$ws.Range("O9").Value = 0.2090378408570124
$ws.Range("P9").Value = 0.2090378408570124
$ws.Range("Q9").Value = 208.9285136243189
$ws.Range("R9").Value = 1880.35662261887
$ws.Range("S9").Value = 0.05660346354158893
$ws.Range("T9").Value = 0.05660346354158893

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 44.17507666666666
$ws.Range("H10").Value = 132.52523
$ws.Range("I10").Value = 0.2707809423859638
$ws.Range("J10").Value = 0.2707809423859638
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.783812
$ws.Range("N10").Value = 44.351436
$ws.Range("O10").Value = 0.6534177673993219
$ws.Range("P10").Value = 0.6534177673993219
$ws.Range("Q10").Value = 653.0760285255866
$ws.Range("R10").Value = 653.0760285255866
$ws.Range("S10").Value = 0.1769330788281209
$ws.Range("T10").Value = 0.1769330788281209
